$d = $word.ActiveDocument

# Step 1: remove the comma that currently trails "oys," (inside the <m>boys,</m> run)
$d.Content.Find.Execute("oys,", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "oys", 2)

# Step 2: locate the "</m>" that now immediately follows that "oys" and append
# a comma after it, restoring the comma on the far side of the closing tag.
$rng = $d.Content
$rng.Find.Execute("oys</m>", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0)
$ins = $d.Range($rng.End, $rng.End)
$ins.InsertAfter(",")
